$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new Wins/Losses/Ties columns, matching the
# existing header style (bold, centered, bordered) by copying the format
# from an existing header cell (A1) before setting the text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD42").Value = 75
$ws.Range("AE2:AE42").Value = 87
$ws.Range("AF2:AF42").Value = 0
